$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AgazatiBesorolas")
$ws.Select()

$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

$ws.Range("A5").Value = "5N-07"
$ws.Range("B5").Value = "Kasza Elemér"
$ws.Range("C5").Value = "DEFGHI"
$ws.Range("D5").Value = "Aktív"
$ws.Range("E5").Value = "Rendszertervezés"
$ws.Range("F5").Value = "Rendszertervezés"

$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

$ws.Range("B5").Select()
